# The source commit ("fix outfall bug, update documentation") only
# touched other files in the repository. For this workbook
# (test_data/swmm_data/gisswmm_quality.xlsx), the associated OOXML diff
# is purely a re-serialization of existing XML (attributes of
# <workbookView>, <xf>, <cellStyle>, <tableStyles> and <pageMargins>
# written in a different order) with every attribute/value pair left
# unchanged, and no cell values, formulas, styles, sheets or structure
# were added, removed or modified.
#
# There is therefore no content-level change to apply through the
# Excel object model: we simply touch the already-open workbook
# without mutating anything, so that it is resaved unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
